$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 12.15
$ws.Range("D2").Value = 1.899999999999995
$ws.Range("E2").Value = -9.800000000000001
$ws.Range("F2").Value = 20.7
$ws.Range("G2").Value = 16.1
$ws.Range("H2").Value = -8.9

# Row 3
$ws.Range("D3").Value = -6.700000000000005
$ws.Range("E3").Value = 2.1
$ws.Range("F3").Value = 25.7
$ws.Range("G3").Value = 16.1
$ws.Range("H3").Value = 15.4

# Row 4
$ws.Range("D4").Value = 1.899999999999995
$ws.Range("E4").Value = -9.800000000000001
$ws.Range("F4").Value = 20.7
$ws.Range("G4").Value = 16.1
$ws.Range("H4").Value = -8.9

# Row 5
$ws.Range("D5").Value = -6.700000000000005
$ws.Range("E5").Value = 2.1
$ws.Range("F5").Value = 25.7
$ws.Range("G5").Value = 16.1
$ws.Range("H5").Value = 15.4

# Row 6
$ws.Range("D6").Value = 1.899999999999995
$ws.Range("E6").Value = -9.800000000000001
$ws.Range("F6").Value = 20.7
$ws.Range("G6").Value = 16.1
$ws.Range("H6").Value = -8.9

# Row 7
$ws.Range("D7").Value = -6.700000000000005
$ws.Range("E7").Value = 2.1
$ws.Range("F7").Value = 25.7
$ws.Range("G7").Value = 16.1
$ws.Range("H7").Value = 15.4

# Row 8
$ws.Range("D8").Value = -6.700000000000005
$ws.Range("E8").Value = 2.1
$ws.Range("F8").Value = 25.7
$ws.Range("G8").Value = 16.1
$ws.Range("H8").Value = 15.4

# Row 9
$ws.Range("D9").Value = 1.899999999999995
$ws.Range("E9").Value = -9.800000000000001
$ws.Range("F9").Value = 20.7
$ws.Range("G9").Value = 16.1
$ws.Range("H9").Value = -8.9

# Row 10
$ws.Range("D10").Value = -6.700000000000005
$ws.Range("E10").Value = 2.1
$ws.Range("F10").Value = 25.7
$ws.Range("G10").Value = 16.1
$ws.Range("H10").Value = 15.4

# Row 11
$ws.Range("D11").Value = 1.899999999999995
$ws.Range("E11").Value = -9.800000000000001
$ws.Range("F11").Value = 20.7
$ws.Range("G11").Value = 16.1
$ws.Range("H11").Value = -8.9

# Row 12
$ws.Range("D12").Value = -6.700000000000005
$ws.Range("E12").Value = 2.1
$ws.Range("F12").Value = 25.7
$ws.Range("G12").Value = 16.1
$ws.Range("H12").Value = 15.4
